$d = $word.ActiveDocument

$replacements = @(
    @("28×17=", "39×79="),
    @("16×28=", "14×38="),
    @("98×79=", "59×35="),
    @("18×77=", "79×91="),
    @("29×32=", "60×92="),
    @("87×98=", "56×98="),
    @("49×53=", "74×43="),
    @("56×73=", "78×93="),
    @("94×12=", "57×58="),
    @("28×84=", "74×97="),
    @("60×39=", "66×69="),
    @("41×37=", "57×41="),
    @("30×51=", "80×75="),
    @("63×43=", "13×44="),
    @("43×84=", "82×41="),
    @("27×69=", "88×46="),
    @("46×78=", "43×87="),
    @("32×97=", "69×87="),
    @("23×69=", "55×94="),
    @("50×45=", "33×30="),
    @("56×82=", "51×81="),
    @("30×40=", "84×11="),
    @("54×23=", "17×45="),
    @("34×40=", "37×17="),
    @("17×71=", "42×70=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
